# feat: add 2022-Q1 data
#
# Before: sheets are 2021-Q2, 2021-Q3, 2021-Q4, 总计 (a rolling summary).
# After:  sheets are 2021-Q2, 2021-Q3, 2021-Q4, 2022-Q1, 总计 - a new
#         "2022-Q1" fund-holdings detail sheet (same shape as the other
#         quarters) is inserted right before the summary sheet, and the
#         summary sheet gets a new leading row for the 2022-Q1 totals.
#
# To reproduce the sheetId numbering of the real commit (2022-Q1 ends up
# with the *old* 总计's sheetId, and the rebuilt 总计 gets a fresh one),
# we repurpose the existing "总计" worksheet object as "2022-Q1" and add
# a brand-new worksheet named 总计 right after it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: grab the header/index-column formatting (bold + thin border +
# center/top alignment, style index "2" in the source file) from the
# existing "总计" sheet *before* we touch any values, so we can stamp it
# onto the newly-needed cells (E1:H1 and A5:A16) without inventing a new
# style definition.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")

$q1.Range("B1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$q1.Range("A2").Copy()
$q1.Range("A2:A16").PasteSpecial(-4122)

$q1.Name = "2022-Q1"

# Header row
$q1.Cells.Item(1, 2).Value = "基金代码"
$q1.Cells.Item(1, 3).Value = "基金名称"
$q1.Cells.Item(1, 4).Value = "基金规模"
$q1.Cells.Item(1, 5).Value = "股票总仓位"
$q1.Cells.Item(1, 6).Value = "仓位占比"
$q1.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q1.Cells.Item(1, 8).Value = "仓位排名"

# Fund holding rows: code, name, fund scale, total stock position,
# position ratio, held market value (100M yuan), position rank.
$fundRows = @(
    @("010610", "上投摩根远见两年持有期混合",     "56.89", "88.51", "3.60", "2.0480", 8),
    @("375010", "上投摩根中国优势混合",           "21.68", "86.19", "3.28", "0.7111", 8),
    @("011046", "富国优质企业混合A",               "8.18",  "71.23", "4.61", "0.3771", 3),
    @("009782", "富国兴泉回报12个月持有期混合A",   "6.18",  "70.06", "4.69", "0.2898", 3),
    @("010029", "富国稳进回报12个月持有期混合A",   "9.16",  "30.21", "1.46", "0.1337", 6),
    @("013678", "富国信享回报12个月持有期混合A",   "9.49",  "27.59", "1.34", "0.1272", 6),
    @("009783", "富国兴泉回报12个月持有期混合C",   "2.45",  "70.06", "4.69", "0.1149", 3),
    @("005732", "富国臻选成长灵活配置混合",         "2.45",  "64.81", "3.96", "0.0970", 4),
    @("005593", "上投摩根创新商业模式灵活配置混合", "1.95",  "91.28", "3.01", "0.0587", 8),
    @("013679", "富国信享回报12个月持有期混合C",   "2.76",  "27.59", "1.34", "0.0370", 6),
    @("005459", "银河嘉谊灵活配置混合A",           "6.47",  "39.69", "0.54", "0.0349", 8),
    @("011047", "富国优质企业混合C",               "0.48",  "71.23", "4.61", "0.0221", 3),
    @("004557", "北信瑞丰鼎丰灵活配置混合",         "0.39",  "64.13", "5.10", "0.0199", 6),
    @("010030", "富国稳进回报12个月持有期混合C",   "1.29",  "30.21", "1.46", "0.0188", 6),
    @("005460", "银河嘉谊灵活配置混合C",           "2.79",  "39.69", "0.54", "0.0151", 8)
)

$row = 2
foreach ($fund in $fundRows) {
    $q1.Cells.Item($row, 1).Value = ($row - 2)
    # fund code / scale / positions look numeric -> force text so leading
    # zeros and trailing zeros survive, same as the source workbook.
    $q1.Cells.Item($row, 2).Value = "'" + $fund[0]
    $q1.Cells.Item($row, 3).Value = $fund[1]
    $q1.Cells.Item($row, 4).Value = "'" + $fund[2]
    $q1.Cells.Item($row, 5).Value = "'" + $fund[3]
    $q1.Cells.Item($row, 6).Value = "'" + $fund[4]
    $q1.Cells.Item($row, 7).Value = "'" + $fund[5]
    $q1.Cells.Item($row, 8).Value = $fund[6]
    $row++
}

# ---------------------------------------------------------------------
# Step 2: create a brand new "总计" summary sheet right after 2022-Q1
# (this gets a fresh sheetId, matching the target commit).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

# Stamp the same header/index-column style onto the new sheet, copying
# from the now-finalized "2022-Q1" sheet.
$q1.Range("B1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)

$q1.Range("A2").Copy()
$total.Range("A2:A5").PasteSpecial(-4122)

$total.Cells.Item(1, 2).Value = "日期"
$total.Cells.Item(1, 3).Value = "持有数量(只)"
$total.Cells.Item(1, 4).Value = "持有市值(亿元)"

$summaryRows = @(
    @("2022-Q1", 15, 4.11),
    @("2021-Q4", 26, 11.35),
    @("2021-Q3", 24, 5.32),
    @("2021-Q2", 3, 0.08)
)

$row = 2
foreach ($entry in $summaryRows) {
    $total.Cells.Item($row, 1).Value = ($row - 2)
    $total.Cells.Item($row, 2).Value = $entry[0]
    $total.Cells.Item($row, 3).Value = $entry[1]
    $total.Cells.Item($row, 4).Value = $entry[2]
    $row++
}

$total.Range("A1").Select()
